# Weekly refresh of the Damasco (Terminal Hortofruticola Agro Chillan) price table:
# existing rows 2-19 are updated in place with the new weekly figures and a brand
# new observation is appended as row 20 (sheet dimension grows from A1:T19 to A1:T20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 45264  # D2 (Fecha)
$ws.Cells.Item(2, 11).Value = 'Dina'  # K2 (Variedad)
$ws.Cells.Item(2, 13).Value = 50  # M2 (Volumen)
$ws.Cells.Item(2, 14).Value = 24000  # N2 (Precio mínimo)
$ws.Cells.Item(2, 15).Value = 24000  # O2 (Precio máximo)
$ws.Cells.Item(2, 16).Value = 24000  # P2 (Precio promedio ponderado)
$ws.Cells.Item(2, 17).Value = '$/caja 10 kilos'  # Q2 (Unidad de comercialización)
$ws.Cells.Item(2, 18).Value = 'Provincia de Limarí'  # R2 (Origen)
$ws.Cells.Item(2, 19).Value = 2400  # S2 (Precio $/Kg)

# Row 3
$ws.Cells.Item(3, 4).Value = 44544  # D3 (Fecha)
$ws.Cells.Item(3, 13).Value = 160  # M3 (Volumen)
$ws.Cells.Item(3, 14).Value = 16000  # N3 (Precio mínimo)
$ws.Cells.Item(3, 15).Value = 17000  # O3 (Precio máximo)
$ws.Cells.Item(3, 16).Value = 16500  # P3 (Precio promedio ponderado)
$ws.Cells.Item(3, 17).Value = '$/caja 15 kilos'  # Q3 (Unidad de comercialización)
$ws.Cells.Item(3, 19).Value = 1100  # S3 (Precio $/Kg)
$ws.Cells.Item(3, 20).Value = 15  # T3 (Kg / unidad)

# Row 4
$ws.Cells.Item(4, 4).Value = 44181  # D4 (Fecha)
$ws.Cells.Item(4, 11).Value = 'Modesto'  # K4 (Variedad)
$ws.Cells.Item(4, 13).Value = 50  # M4 (Volumen)
$ws.Cells.Item(4, 14).Value = 20000  # N4 (Precio mínimo)
$ws.Cells.Item(4, 15).Value = 21000  # O4 (Precio máximo)
$ws.Cells.Item(4, 16).Value = 20500  # P4 (Precio promedio ponderado)
$ws.Cells.Item(4, 18).Value = 'Región de Coquimbo'  # R4 (Origen)
$ws.Cells.Item(4, 19).Value = 1139  # S4 (Precio $/Kg)

# Row 5
$ws.Cells.Item(5, 4).Value = 44176  # D5 (Fecha)
$ws.Cells.Item(5, 13).Value = 50  # M5 (Volumen)
$ws.Cells.Item(5, 14).Value = 17000  # N5 (Precio mínimo)
$ws.Cells.Item(5, 15).Value = 18000  # O5 (Precio máximo)
$ws.Cells.Item(5, 16).Value = 17400  # P5 (Precio promedio ponderado)
$ws.Cells.Item(5, 17).Value = '$/caja 18 kilos'  # Q5 (Unidad de comercialización)
$ws.Cells.Item(5, 19).Value = 967  # S5 (Precio $/Kg)
$ws.Cells.Item(5, 20).Value = 18  # T5 (Kg / unidad)

# Row 6
$ws.Cells.Item(6, 4).Value = 45267  # D6 (Fecha)
$ws.Cells.Item(6, 11).Value = 'Dina'  # K6 (Variedad)
$ws.Cells.Item(6, 13).Value = 100  # M6 (Volumen)
$ws.Cells.Item(6, 14).Value = 22000  # N6 (Precio mínimo)
$ws.Cells.Item(6, 15).Value = 22000  # O6 (Precio máximo)
$ws.Cells.Item(6, 16).Value = 22000  # P6 (Precio promedio ponderado)
$ws.Cells.Item(6, 17).Value = '$/caja 10 kilos'  # Q6 (Unidad de comercialización)
$ws.Cells.Item(6, 18).Value = 'Provincia de Limarí'  # R6 (Origen)
$ws.Cells.Item(6, 19).Value = 2200  # S6 (Precio $/Kg)
$ws.Cells.Item(6, 20).Value = 10  # T6 (Kg / unidad)

# Row 7
$ws.Cells.Item(7, 4).Value = 44187  # D7 (Fecha)
$ws.Cells.Item(7, 11).Value = 'Dina'  # K7 (Variedad)
$ws.Cells.Item(7, 12).Value = 'Primera'  # L7 (Calidad)
$ws.Cells.Item(7, 13).Value = 55  # M7 (Volumen)
$ws.Cells.Item(7, 14).Value = 15000  # N7 (Precio mínimo)
$ws.Cells.Item(7, 15).Value = 16000  # O7 (Precio máximo)
$ws.Cells.Item(7, 16).Value = 15455  # P7 (Precio promedio ponderado)
$ws.Cells.Item(7, 17).Value = '$/caja 15 kilos granel'  # Q7 (Unidad de comercialización)
$ws.Cells.Item(7, 19).Value = 1030  # S7 (Precio $/Kg)

# Row 8
$ws.Cells.Item(8, 4).Value = 44165  # D8 (Fecha)
$ws.Cells.Item(8, 12).Value = 'Segunda'  # L8 (Calidad)
$ws.Cells.Item(8, 14).Value = 16000  # N8 (Precio mínimo)
$ws.Cells.Item(8, 15).Value = 17000  # O8 (Precio máximo)
$ws.Cells.Item(8, 16).Value = 16500  # P8 (Precio promedio ponderado)
$ws.Cells.Item(8, 17).Value = '$/caja 15 kilos granel'  # Q8 (Unidad de comercialización)
$ws.Cells.Item(8, 18).Value = 'Provincia de Limarí'  # R8 (Origen)
$ws.Cells.Item(8, 19).Value = 1100  # S8 (Precio $/Kg)
$ws.Cells.Item(8, 20).Value = 15  # T8 (Kg / unidad)

# Row 9
$ws.Cells.Item(9, 4).Value = 44551  # D9 (Fecha)
$ws.Cells.Item(9, 12).Value = 'Primera'  # L9 (Calidad)
$ws.Cells.Item(9, 13).Value = 120  # M9 (Volumen)
$ws.Cells.Item(9, 14).Value = 15500  # N9 (Precio mínimo)
$ws.Cells.Item(9, 15).Value = 16000  # O9 (Precio máximo)
$ws.Cells.Item(9, 16).Value = 15750  # P9 (Precio promedio ponderado)
$ws.Cells.Item(9, 17).Value = '$/caja 15 kilos'  # Q9 (Unidad de comercialización)
$ws.Cells.Item(9, 19).Value = 1050  # S9 (Precio $/Kg)
$ws.Cells.Item(9, 20).Value = 15  # T9 (Kg / unidad)

# Row 10
$ws.Cells.Item(10, 4).Value = 44907  # D10 (Fecha)
$ws.Cells.Item(10, 11).Value = 'Castle Brite'  # K10 (Variedad)
$ws.Cells.Item(10, 13).Value = 120  # M10 (Volumen)
$ws.Cells.Item(10, 14).Value = 15000  # N10 (Precio mínimo)
$ws.Cells.Item(10, 15).Value = 16000  # O10 (Precio máximo)
$ws.Cells.Item(10, 16).Value = 15500  # P10 (Precio promedio ponderado)
$ws.Cells.Item(10, 17).Value = '$/bandeja 10 kilos'  # Q10 (Unidad de comercialización)
$ws.Cells.Item(10, 18).Value = 'Región de O''Higgins'  # R10 (Origen)
$ws.Cells.Item(10, 19).Value = 1550  # S10 (Precio $/Kg)

# Row 11
$ws.Cells.Item(11, 4).Value = 44907  # D11 (Fecha)
$ws.Cells.Item(11, 12).Value = 'Segunda'  # L11 (Calidad)
$ws.Cells.Item(11, 13).Value = 60  # M11 (Volumen)
$ws.Cells.Item(11, 14).Value = 14000  # N11 (Precio mínimo)
$ws.Cells.Item(11, 15).Value = 14000  # O11 (Precio máximo)
$ws.Cells.Item(11, 16).Value = 14000  # P11 (Precio promedio ponderado)
$ws.Cells.Item(11, 17).Value = '$/bandeja 10 kilos'  # Q11 (Unidad de comercialización)
$ws.Cells.Item(11, 19).Value = 1400  # S11 (Precio $/Kg)
$ws.Cells.Item(11, 20).Value = 10  # T11 (Kg / unidad)

# Row 12
$ws.Cells.Item(12, 4).Value = 45273  # D12 (Fecha)
$ws.Cells.Item(12, 13).Value = 80  # M12 (Volumen)
$ws.Cells.Item(12, 14).Value = 20000  # N12 (Precio mínimo)
$ws.Cells.Item(12, 15).Value = 20000  # O12 (Precio máximo)
$ws.Cells.Item(12, 16).Value = 20000  # P12 (Precio promedio ponderado)
$ws.Cells.Item(12, 17).Value = '$/bandeja 10 kilos'  # Q12 (Unidad de comercialización)
$ws.Cells.Item(12, 19).Value = 2000  # S12 (Precio $/Kg)
$ws.Cells.Item(12, 20).Value = 10  # T12 (Kg / unidad)

# Row 13
$ws.Cells.Item(13, 4).Value = 44904  # D13 (Fecha)
$ws.Cells.Item(13, 11).Value = 'Castle Brite'  # K13 (Variedad)
$ws.Cells.Item(13, 13).Value = 60  # M13 (Volumen)
$ws.Cells.Item(13, 14).Value = 15000  # N13 (Precio mínimo)
$ws.Cells.Item(13, 15).Value = 16000  # O13 (Precio máximo)
$ws.Cells.Item(13, 16).Value = 15500  # P13 (Precio promedio ponderado)
$ws.Cells.Item(13, 17).Value = '$/bandeja 10 kilos'  # Q13 (Unidad de comercialización)
$ws.Cells.Item(13, 18).Value = 'Región de O''Higgins'  # R13 (Origen)
$ws.Cells.Item(13, 19).Value = 1550  # S13 (Precio $/Kg)

# Row 14
$ws.Cells.Item(14, 4).Value = 44904  # D14 (Fecha)
$ws.Cells.Item(14, 12).Value = 'Segunda'  # L14 (Calidad)
$ws.Cells.Item(14, 14).Value = 14000  # N14 (Precio mínimo)
$ws.Cells.Item(14, 15).Value = 14000  # O14 (Precio máximo)
$ws.Cells.Item(14, 16).Value = 14000  # P14 (Precio promedio ponderado)
$ws.Cells.Item(14, 17).Value = '$/bandeja 10 kilos'  # Q14 (Unidad de comercialización)
$ws.Cells.Item(14, 18).Value = 'Región de O''Higgins'  # R14 (Origen)
$ws.Cells.Item(14, 19).Value = 1400  # S14 (Precio $/Kg)
$ws.Cells.Item(14, 20).Value = 10  # T14 (Kg / unidad)

# Row 15
$ws.Cells.Item(15, 4).Value = 44552  # D15 (Fecha)
$ws.Cells.Item(15, 13).Value = 120  # M15 (Volumen)
$ws.Cells.Item(15, 14).Value = 15500  # N15 (Precio mínimo)
$ws.Cells.Item(15, 15).Value = 16000  # O15 (Precio máximo)
$ws.Cells.Item(15, 16).Value = 15750  # P15 (Precio promedio ponderado)
$ws.Cells.Item(15, 17).Value = '$/caja 15 kilos'  # Q15 (Unidad de comercialización)
$ws.Cells.Item(15, 19).Value = 1050  # S15 (Precio $/Kg)
$ws.Cells.Item(15, 20).Value = 15  # T15 (Kg / unidad)

# Row 16
$ws.Cells.Item(16, 4).Value = 44537  # D16 (Fecha)
$ws.Cells.Item(16, 12).Value = 'Primera'  # L16 (Calidad)
$ws.Cells.Item(16, 14).Value = 21000  # N16 (Precio mínimo)
$ws.Cells.Item(16, 15).Value = 21500  # O16 (Precio máximo)
$ws.Cells.Item(16, 16).Value = 21250  # P16 (Precio promedio ponderado)
$ws.Cells.Item(16, 17).Value = '$/caja 15 kilos'  # Q16 (Unidad de comercialización)
$ws.Cells.Item(16, 18).Value = 'Región de O''Higgins'  # R16 (Origen)
$ws.Cells.Item(16, 19).Value = 1417  # S16 (Precio $/Kg)

# Row 17
$ws.Cells.Item(17, 4).Value = 44189  # D17 (Fecha)
$ws.Cells.Item(17, 11).Value = 'Dina'  # K17 (Variedad)
$ws.Cells.Item(17, 13).Value = 80  # M17 (Volumen)
$ws.Cells.Item(17, 14).Value = 16000  # N17 (Precio mínimo)
$ws.Cells.Item(17, 15).Value = 17000  # O17 (Precio máximo)
$ws.Cells.Item(17, 16).Value = 16562  # P17 (Precio promedio ponderado)
$ws.Cells.Item(17, 18).Value = 'Región de O''Higgins'  # R17 (Origen)
$ws.Cells.Item(17, 19).Value = 920  # S17 (Precio $/Kg)

# Row 18
$ws.Cells.Item(18, 4).Value = 44168  # D18 (Fecha)
$ws.Cells.Item(18, 11).Value = 'Castle Brite'  # K18 (Variedad)
$ws.Cells.Item(18, 13).Value = 30  # M18 (Volumen)
$ws.Cells.Item(18, 14).Value = 16000  # N18 (Precio mínimo)
$ws.Cells.Item(18, 15).Value = 17000  # O18 (Precio máximo)
$ws.Cells.Item(18, 16).Value = 16500  # P18 (Precio promedio ponderado)
$ws.Cells.Item(18, 17).Value = '$/caja 16 kilos granel'  # Q18 (Unidad de comercialización)
$ws.Cells.Item(18, 18).Value = 'Región de Coquimbo'  # R18 (Origen)
$ws.Cells.Item(18, 19).Value = 1031  # S18 (Precio $/Kg)
$ws.Cells.Item(18, 20).Value = 16  # T18 (Kg / unidad)

# Row 19
$ws.Cells.Item(19, 4).Value = 44174  # D19 (Fecha)
$ws.Cells.Item(19, 11).Value = 'Castle Brite'  # K19 (Variedad)
$ws.Cells.Item(19, 13).Value = 75  # M19 (Volumen)
$ws.Cells.Item(19, 14).Value = 9000  # N19 (Precio mínimo)
$ws.Cells.Item(19, 15).Value = 10000  # O19 (Precio máximo)
$ws.Cells.Item(19, 16).Value = 9467  # P19 (Precio promedio ponderado)
$ws.Cells.Item(19, 17).Value = '$/caja 10 kilos'  # Q19 (Unidad de comercialización)
$ws.Cells.Item(19, 19).Value = 947  # S19 (Precio $/Kg)
$ws.Cells.Item(19, 20).Value = 10  # T19 (Kg / unidad)

# Row 20
$ws.Cells.Item(20, 1).Value = 7  # A20 (Mercado ID)
$ws.Cells.Item(20, 2).Value = 'Terminal Hortofrutícola Agro Chillán'  # B20 (Mercado)
$ws.Cells.Item(20, 3).Value = 'Ñuble'  # C20 (Región)
$ws.Cells.Item(20, 4).Value = 45265  # D20 (Fecha)
$ws.Cells.Item(20, 5).Value = 16  # E20 (Codreg)
$ws.Cells.Item(20, 6).Value = 'Fruta'  # F20 (Tipo)
$ws.Cells.Item(20, 7).Value = 100103  # G20 (Producto ID)
$ws.Cells.Item(20, 8).Value = 'Frutos de hueso (carozo)'  # H20 (Producto)
$ws.Cells.Item(20, 9).Value = 100103003  # I20 (Categoría ID)
$ws.Cells.Item(20, 10).Value = 'Damasco'  # J20 (Categoría)
$ws.Cells.Item(20, 11).Value = 'Dina'  # K20 (Variedad)
$ws.Cells.Item(20, 12).Value = 'Primera'  # L20 (Calidad)
$ws.Cells.Item(20, 13).Value = 80  # M20 (Volumen)
$ws.Cells.Item(20, 14).Value = 20000  # N20 (Precio mínimo)
$ws.Cells.Item(20, 15).Value = 20000  # O20 (Precio máximo)
$ws.Cells.Item(20, 16).Value = 20000  # P20 (Precio promedio ponderado)
$ws.Cells.Item(20, 17).Value = '$/caja 10 kilos'  # Q20 (Unidad de comercialización)
$ws.Cells.Item(20, 18).Value = 'Provincia de Limarí'  # R20 (Origen)
$ws.Cells.Item(20, 19).Value = 2000  # S20 (Precio $/Kg)
$ws.Cells.Item(20, 20).Value = 10  # T20 (Kg / unidad)

# Give the new row's Fecha cell (D20) the same date number format used by the
# rest of column D (copied from D3, an existing/unaffected date cell).
$ws.Cells.Item(20, 4).NumberFormat = $ws.Cells.Item(3, 4).NumberFormat
